$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 197, pushing existing rows 197:205 down to 198:206.
$ws.Rows(197).Insert()

# Populate the newly inserted row 197 with the new weekly record.
$ws.Range("A197").Value2 = 5
$ws.Range("B197").Value2 = "Macroferia Regional de Talca"
$ws.Range("C197").Value2 = "Maule"
$ws.Range("D197").Value2 = 44931
$ws.Range("E197").Value2 = 7
$ws.Range("F197").Value2 = 100112031
$ws.Range("G197").Value2 = "Poroto verde"
$ws.Range("H197").Value2 = "Sin especificar"
$ws.Range("I197").Value2 = "Primera"
$ws.Range("J197").Value2 = 150
$ws.Range("K197").Value2 = 35000
$ws.Range("L197").Value2 = 35000
$ws.Range("M197").Value2 = 35000
$ws.Range("N197").Value2 = "$/saco 25 kilos"
$ws.Range("O197").Value2 = "Región del Maule"
$ws.Range("P197").Value2 = 1400
$ws.Range("Q197").Value2 = 25
$ws.Range("R197").Value2 = "Hortaliza"
